$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 6353
$ws.Range("C24").Value = 998
$ws.Range("D24").Value = 5955596
$ws.Range("E24").Value = 937.4462458680938
$ws.Range("F24").Value = 8.302079781793381
$ws.Range("G24").Value = 3.419689119170988
$ws.Range("H24").Value = 26.16259216044157
